$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The existing "_GoBack" bookmark (left over from the author's last
#    editing session) sat at the end of the "Variants" paragraph. Since
#    we are about to make a new edit elsewhere, that stale bookmark is
#    removed first ...
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Remove the stray leading hyphen from the "Rotate:" bullet.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "-Baby shrew will rotate back and forth before picking an angle at which to move at. The wheels will rotate in opposite directions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Baby shrew will rotate back and forth before picking an angle at which to move at. The wheels will rotate in opposite directions.",
    2) | Out-Null

# ------------------------------------------------------------------
# 3. Remove the stray leading hyphen from the "Death:" bullet.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "-The shrew will fall on its side before exploding.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The shrew will fall on its side before exploding.",
    2) | Out-Null

# ------------------------------------------------------------------
# 4. ... and re-inserted immediately before that last-edited text,
#    since Word always re-drops "_GoBack" at the site of the most
#    recent edit.
# ------------------------------------------------------------------
$targetRange = $d.Content
$targetRange.Find.Execute("The shrew will fall on its side before exploding.") | Out-Null
$newBookmarkRange = $d.Range($targetRange.Start, $targetRange.Start)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
